$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Authorized USBs")

# Update existing "Device E" row's VID_PID value (was 05ac_12a0) to a new
# descriptor value, and add a new entry for "Device F".
$ws.Range("B6").Value = "AAAA_BBBB"
$ws.Range("A7").Value = "Device F"
$ws.Range("B7").Value = "FFFF_2222"

# Leave the newly edited cell selected, matching the saved selection state.
$ws.Range("B7").Select()
